$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RegistrationStatus (column H) for officer rows 2 and 3: "Not Registered" -> "Approved"
$ws.Range("H2").Value = "Approved"
$ws.Range("H3").Value = "Approved"

# Officer row 4 was not registered at all before; now explicitly "Not Registered"
$ws.Range("H4").Value = "Not Registered"

# Set RegisteredProject (column I) for the now-approved officers
$ws.Range("I2").Value = "Acacia Breeze"
$ws.Range("I3").Value = "Acacia Breeze"

# Move active selection to I4 (matches the saved selection state in the file)
$ws.Range("I4").Select()
